$d = $word.ActiveDocument

$rng = $d.Content
$rng.Start = 0
$rng.End = 0

$found = $true
while ($found) {
    $found = $rng.Find.Execute("^t", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.NameAscii = "Helvetica"
        $rng.Font.NameFarEast = "Helvetica"
        $rng.Font.NameOther = "Helvetica"
        $rng.Font.NameBi = "Helvetica"
        $rng.Font.Size = 12
        $rng.Collapse(0)
    }
}
